$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.180.17"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "3.463.16"
$ws.Range("E3").Value = "  +2.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.50%  "

$ws.Range("D7").Value = "3.462.79"
$ws.Range("E7").Value = "  +2.48%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.67"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").Value = "4.058.96"
$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.75"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.03%  "

$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").Value = "3.465.98"
$ws.Range("E17").Value = "  +2.46%  "

$ws.Range("D18").Value = "62.243.86"
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.61"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.58%  "

$ws.Range("E23").Value = "  +1.72%  "

$ws.Range("D24").Value = "3.597.52"
$ws.Range("E24").Value = "  +2.10%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "72.64"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.89%  "

$ws.Range("E28").Value = "  +1.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.179"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +9.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.86"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.58"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -12.38%  "

$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.30"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("E34").Value = "  +1.32%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.11"
$ws.Range("D36").ClearFormats()

$ws.Range("E37").Value = "  +2.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.08"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.59"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.74"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0800"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.99%  "

$ws.Range("E42").Value = "  +2.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.86"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.19"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.55%  "

$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.21"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").Value = "2.662.27"
$ws.Range("E49").Value = "  +11.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.88"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.55%  "

$ws.Range("E51").Value = "  +0.72%  "
